$d = $word.ActiveDocument

# 1) "No.  76" -> "No.  77"  (order number)
$d.Content.Find.Execute("76", $true, $false, $false, $false, $false, $true, 1, $false, "77", 2) | Out-Null

# 2) CONCEPTO cell: fill in the first three paragraphs (Camisa / Camiseta para vender / Camiseta edison)
#    and VALOR cell: fill first three paragraphs (10000 / 10000 / 4000)
$idx = 0
foreach ($par in $d.Paragraphs) {
    $idx++
    switch ($idx) {
        24 { $par.Range.Text = "Camisa- (asdf)" }
        25 { $par.Range.InsertBefore("Camiseta para vender- (n89d)") }
        26 { $par.Range.InsertBefore("Camiseta  edison- (fasdfa)") }
        39 { $par.Range.Text = "10000" }
        40 { $par.Range.InsertBefore("10000") }
        41 { $par.Range.InsertBefore("4000") }
        58 {
            $rng = $par.Range
            $rng.Find.Execute("6000", $true, $false, $false, $false, $false, $true, 1, $false, "10000", 2) | Out-Null
        }
        61 { $par.Range.Text = "2018-11-12" }
        63 {
            $rng = $par.Range
            $rng.Find.Execute("44000", $true, $false, $false, $false, $false, $true, 1, $false, "14000", 2) | Out-Null
        }
        64 { $par.Range.Text = "24000" }
        67 {
            $rng = $par.Range
            $rng.Find.Execute(" - Tarjeta profesional", $true, $false, $false, $false, $false, $true, 1, $false, " - ", 2) | Out-Null
        }
    }
}
